$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Digital Image Processing with Python"
$ws.Range("B2").Value = "Spring 2026"
$ws.Range("C2").Value = "IIT Goa"

$ws.Range("A3").Value = "Deep Learning for Computer Vision"
$ws.Range("B3").Value = "Autumn 2025"
$ws.Range("C3").Value = "IIT Goa"

$ws.Range("A4").Value = "Ocean Vision AI"
$ws.Range("B4").Value = "Autumn 2025"
$ws.Range("C4").Value = "IIT Goa"

$ws.Columns.Item(1).ColumnWidth = 32.5
$ws.Columns.Item(2).ColumnWidth = 11.67

$ws.Range("H19").Select()
